# Update "想去人数" (number of people interested) values in column F
# across the "展览" (sheet1), "演出" (sheet2) and "全部类型" (sheet4) sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 1238
$ws1.Range("F7").Value = 59
$ws1.Range("F10").Value = 3449
$ws1.Range("F16").Value = 591
$ws1.Range("F17").Value = 85
$ws1.Range("F18").Value = 728
$ws1.Range("F20").Value = 120
$ws1.Range("F24").Value = 2587
$ws1.Range("F25").Value = 5098
$ws1.Range("F30").Value = 281
$ws1.Range("F31").Value = 2237
$ws1.Range("F34").Value = 82
$ws1.Range("F35").Value = 111
$ws1.Range("F36").Value = 173
$ws1.Range("F38").Value = 459
$ws1.Range("F39").Value = 793
$ws1.Range("F41").Value = 451
$ws1.Range("F43").Value = 472

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 70

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1238
$ws4.Range("F7").Value = 59
$ws4.Range("F10").Value = 3449
$ws4.Range("F14").Value = 70
$ws4.Range("F17").Value = 591
$ws4.Range("F18").Value = 85
$ws4.Range("F19").Value = 728
$ws4.Range("F21").Value = 120
$ws4.Range("F25").Value = 2587
$ws4.Range("F26").Value = 5098
$ws4.Range("F31").Value = 281
$ws4.Range("F32").Value = 2237
$ws4.Range("F35").Value = 82
$ws4.Range("F36").Value = 111
$ws4.Range("F37").Value = 173
$ws4.Range("F39").Value = 459
$ws4.Range("F40").Value = 793
$ws4.Range("F42").Value = 451
$ws4.Range("F44").Value = 472
